# Deploying to gh-pages: extend the table with the 2022 column (S) and
# refresh the 2020/2021 (Q/R) values for the oblast rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New year header in S4, formatted like the existing year headers (R4) ---
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S4").Value = 2022

# --- 2. Updated Q/R figures for the data rows (5-14) ---
$qrValues = @{
    5  = @{ Q = 117.60684979252385; R = 113.34848864817617 }
    6  = @{ Q = 114.77319768114526; R = 115.06069350712495 }
    7  = @{ Q = 116.40044011407315; R = 114.29658549692938 }
    8  = @{ Q = 117.53828537152096; R = 113.75761785228545 }
    9  = @{ Q = 117.42206669681742; R = 113.98264089946031 }
    10 = @{ Q = 113.98326995089161; R = 113.92720567782911 }
    11 = @{ Q = 123.488978736909;   R = 114.17226706705155 }
    12 = @{ Q = 118.12340252754679; R = 114.45153946490467 }
    13 = @{ Q = 118.87059844457349; R = 112.69493421065988 }
    14 = @{ Q = 114.06377070452145; R = 113.95067699644588 }
}

foreach ($row in $qrValues.Keys) {
    $ws.Range("Q$row").Value = $qrValues[$row].Q
    $ws.Range("R$row").Value = $qrValues[$row].R
}

# --- 3. New 2022 column (S) values for rows 5-13, formatted with the plain
#        default style (matches column A's unstyled numeric look, s=4) ---
$ws.Range("A6").Copy()
$ws.Range("S5:S13").PasteSpecial(-4122)   # xlPasteFormats

$sValues = @{
    5  = 115.8
    6  = 115.2
    7  = 115.4
    8  = 111.8
    9  = 116.8
    10 = 108.2
    11 = 111
    12 = 115.8
    13 = 117.9
}
foreach ($row in $sValues.Keys) {
    $ws.Range("S$row").Value = $sValues[$row]
}

# --- 4. Row 14 (totals row) S cell, formatted like R14 (thick-bottom style) ---
$ws.Range("R14").Copy()
$ws.Range("S14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S14").Value = 112.4

# --- 5. Selection moves from T6 to T4 ---
[void]$ws.Range("T4").Select()

Write-Output "edit applied"
